# This script applies the edit described by the commit "Ref til Python fil oppg. B":
#  1. "Kode er: " (bold+underline heading) is re-typed, which Word's proofing
#     engine marks up with <w:proofErr> spell-check ranges around "Kode" and "er".
#  2. The following "B C C A C B C C C C C C C C C C C A" line is likewise
#     re-split into runs with <w:proofErr> markers around the repeated lone "C"
#     tokens (this happens twice in the document -- once under "Kode er:" and
#     once under "d) Huffman kode:").
#  3. The paragraph block consisting of the "Ascii - 144 bits ..." remark, a
#     blank paragraph, "b) " and "c) " is moved up so it immediately follows
#     the "... bits (ascii tabellen ...)" paragraph, and the "b) " paragraph
#     gets a new bold-italic reference appended: "Se Python fil Oppgave 1.2.1 B".
#     The _GoBack bookmark travels with the block (it always sits on the last
#     paragraph of that group).

$d = $word.ActiveDocument

function New-PackageXml([string]$bodyFragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $bodyFragment +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Get-ParaText($para) {
    # Paragraph.Range.Text includes the trailing paragraph-mark character(s);
    # strip them so callers can compare against plain text.
    return $para.Range.Text.TrimEnd([char]13, [char]7)
}

function Find-ParagraphIndex([string]$exactText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ((Get-ParaText($d.Paragraphs($i))) -eq $exactText) {
            return $i
        }
    }
    return -1
}

# Re-usable run sequence for "B C C A C B C C C C C C C C C C C A" with the
# proofErr-wrapped repeated lone "C" runs, matching how Word re-segments the
# text when it is retyped with as-you-type spell checking on.
$bcodeRuns =
    '<w:r><w:t xml:space="preserve">B C </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> A C B C </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>C</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> A</w:t></w:r>'

# ---------------------------------------------------------------------------
# Step 1: Move the "Ascii - 144 bits..."/blank/"b) "/"c) " block so it follows
# the "... bits (ascii tabellen ...)" paragraph, modify "b) " to add the
# Python-file reference, and carry the _GoBack bookmark with the block.
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$asciiIdx = Find-ParagraphIndex("Ascii – 144 bits – hvert symbol 8 bits ( iflg tekst fra uke 5 er hvert symbol 8 bits)")
$cIdx = Find-ParagraphIndex("c) ")

$startP = $d.Paragraphs($asciiIdx)
$endP = $d.Paragraphs($cIdx)
$moveRange = $d.Range($startP.Range.Start, $endP.Range.End)

$movedBlockFragment = (
    '<w:p><w:pPr><w:rPr><w:lang w:val="nb-NO"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="nb-NO"/></w:rPr><w:t>Ascii &#8211; 144 bits &#8211; hvert symbol 8 bits</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="nb-NO"/></w:rPr><w:t xml:space="preserve"> ( iflg tekst fra uke 5</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="nb-NO"/></w:rPr><w:t xml:space="preserve"> er hvert symbol 8 bits</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="nb-NO"/></w:rPr><w:t>)</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="nb-NO"/></w:rPr></w:pPr></w:p>' +
    '<w:p><w:pPr><w:rPr><w:b/><w:i/><w:lang w:val="nb-NO"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="nb-NO"/></w:rPr><w:t xml:space="preserve">b) </w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="nb-NO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:i/><w:lang w:val="nb-NO"/></w:rPr><w:t>Se Python fil Oppgave 1.2.1 B</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="nb-NO"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="nb-NO"/></w:rPr><w:t xml:space="preserve">c) </w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
)

$moveRange.InsertXML((New-PackageXml $movedBlockFragment))

# ---------------------------------------------------------------------------
# Step 2: Re-split "Kode er: " into proofErr-wrapped runs.
# ---------------------------------------------------------------------------

$kodeIdx = Find-ParagraphIndex("Kode er: ")
$kodeP = $d.Paragraphs($kodeIdx)
$kodeRange = $kodeP.Range

$kodeFragment =
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Kode</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>er</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>'
$kodePara = '<w:p><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr>' + $kodeFragment + '</w:p>'

$kodeRange.InsertXML((New-PackageXml $kodePara))

# ---------------------------------------------------------------------------
# Step 3: Re-split both "B C C A C B C C C C C C C C C C C A" paragraphs into
# proofErr-wrapped runs (one follows "Kode er: ", the other follows
# "d) Huffman kode:").
# ---------------------------------------------------------------------------

$bcodeIndices = @()
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ((Get-ParaText($d.Paragraphs($i))) -eq "B C C A C B C C C C C C C C C C C A") {
        $bcodeIndices += $i
    }
}

foreach ($idx in $bcodeIndices) {
    $p = $d.Paragraphs($idx)
    $r = $p.Range
    $para = '<w:p>' + $bcodeRuns + '</w:p>'
    $r.InsertXML((New-PackageXml $para))
}

Write-Host "Edit complete"
